$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend column-A formatting (style) from row 2 down through the new rows
# before writing values, so A3:A6 keep the same "s=1" style as A2/A1.
$ws.Range("A2").Copy()
$ws.Range("A3:A6").PasteSpecial(-4122)

# Row 2 - Darul Aman Kabul (2), Afghanistan
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Darul Aman Kabul (2), Afghanistan"
$ws.Range("C2").Value = "Afghanistan"
$ws.Range("D2").Value = "AFG"
$ws.Range("E2").Value = 34.48845
$ws.Range("F2").Value = "Darul Aman Kabul (2)"
$ws.Range("G2").Value = 69.20287999999999
$ws.Range("H2").Value = 34.48845
$ws.Range("I2").Value = 69.20287999999999
$ws.Range("J2").Value = "Asia"
$ws.Range("K2").Value = "correct location data"

# Row 3 - Dinajpur, Bangladesh
$ws.Range("A3").Value = 20
$ws.Range("B3").Value = "Dinajpur, 5216, Rangpur Division, Bangladesh"
$ws.Range("C3").Value = "Bangladesh"
$ws.Range("D3").Value = "BGD"
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = "Dinajpur"
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 25.6260712
$ws.Range("I3").Value = 88.6346228
$ws.Range("J3").Value = "Asia"
$ws.Range("K3").Value = "no lat/lng entered / incorrect lat/lng - geocoded location"

# Row 4 - Bolivia (Address / Recorded_Lat / Recorded_Lng left blank)
$ws.Range("A4").Value = 52
$ws.Range("C4").Value = "Bolivia, Plurinational State of"
$ws.Range("D4").Value = "BOL"
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = "Lomas De Arena (santa Cruz)"
$ws.Range("G4").Value = 0
$ws.Range("J4").Value = "South America"
$ws.Range("K4").Value = "incorrect location data/cannot find coordinates"

# Row 5 - Zambia
$ws.Range("A5").Value = 1024
$ws.Range("B5").Value = "Zamseed, Farm, Zambia"
$ws.Range("C5").Value = "Zambia"
$ws.Range("D5").Value = "ZMB"
$ws.Range("E5").Value = -14.2
$ws.Range("F5").Value = "Zamseed, Farm"
$ws.Range("G5").Value = 28.4
$ws.Range("H5").Value = -14.2
$ws.Range("I5").Value = 28.4
$ws.Range("J5").Value = "Eastern and Southern Africa"
$ws.Range("K5").Value = "correct location data"

# Row 6 - Guadalajara, Mexico
$ws.Range("A6").Value = 599
$ws.Range("B6").Value = "Guadalajara, 44100, Jalisco, Mexico"
$ws.Range("C6").Value = "Mexico"
$ws.Range("D6").Value = "MEX"
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = "Guadalajara, Jalisco"
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 20.676143
$ws.Range("I6").Value = -103.3469982
$ws.Range("J6").Value = "Central America and Caribbean"
$ws.Range("K6").Value = "no lat/lng entered / incorrect lat/lng - geocoded location"
